$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Summary of the change (see commit message "before removing the
# introduction" + the supplied diff):
#   - Paragraph 1 ("מאמר של אמנון כווארי ... – להוסיף ציטוט") moves down to
#     become paragraph 2.
#   - Paragraph 2 ("להוסיף footnote ... בכל ER.") moves down to become
#     paragraph 3.
#   - Paragraph 3 ("לבנות מחדש את הפרק של ה- DYNAMIC ...") is removed.
#   - Paragraph 4 ("לבדוק את הערכים בפרק של ה- DISCUSSION.") is removed
#     (merged away).
#   - Paragraph 1 gets brand-new text: "להעיף את ה- INTRODUCTION (לבדוק
#     שחלקים ממנה מתחילים את הפרקים הרלוונטיים)"
#
# Word's Range/FormattedText objects are *live* (position-bound), so we
# can't simply stash a FormattedText handle for later reuse once earlier
# edits shift character offsets underneath it. Instead we stage copies of
# the two paragraphs that need to move at the very end of the document,
# then paste those staged copies into their final homes (always re-
# resolving Paragraphs(n) fresh immediately before each read/write), and
# finally discard the scratch paragraphs.
# ---------------------------------------------------------------------------

# Append two scratch paragraphs at the end of the document to use as a
# holding area.
$d.Content.InsertParagraphAfter()
$d.Content.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$stageA = $n - 1   # will hold old paragraph 1's text
$stageB = $n        # will hold old paragraph 2's text

# Stage old paragraph 1 ("מאמר של אמנון כווארי ...") into scratch slot A.
$p1 = $d.Paragraphs(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$sa = $d.Paragraphs($stageA)
$ra = $d.Range($sa.Range.Start, $sa.Range.End - 1)
$ra.FormattedText = $r1.FormattedText

# Stage old paragraph 2 ("להוסיף footnote ...") into scratch slot B.
$p2 = $d.Paragraphs(2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$sb = $d.Paragraphs($stageB)
$rb = $d.Range($sb.Range.Start, $sb.Range.End - 1)
$rb.FormattedText = $r2.FormattedText

# Paste staged slot A into paragraph 2 (was "להוסיף footnote ...", now
# becomes "מאמר של אמנון כווארי ...").
$p2now = $d.Paragraphs(2)
$r2now = $d.Range($p2now.Range.Start, $p2now.Range.End - 1)
$saNow = $d.Paragraphs($stageA)
$raNow = $d.Range($saNow.Range.Start, $saNow.Range.End - 1)
$r2now.FormattedText = $raNow.FormattedText

# Paste staged slot B into paragraph 3 (was "לבנות מחדש את הפרק ...", now
# becomes "להוסיף footnote ...").
$p3now = $d.Paragraphs(3)
$r3now = $d.Range($p3now.Range.Start, $p3now.Range.End - 1)
$sbNow = $d.Paragraphs($stageB)
$rbNow = $d.Range($sbNow.Range.Start, $sbNow.Range.End - 1)
$r3now.FormattedText = $rbNow.FormattedText

# Remove paragraph 4 ("לבדוק את הערכים בפרק של ה- DISCUSSION."), which is
# no longer needed (it is effectively merged away in the edit).
$p4now = $d.Paragraphs(4)
$p4now.Range.Delete()

# Remove the two scratch paragraphs (they have shifted down by one since
# paragraph 4 was deleted).
$countNow = $d.Paragraphs.Count
$d.Paragraphs($countNow).Range.Delete()
$d.Paragraphs($countNow - 1).Range.Delete()

# ---------------------------------------------------------------------------
# Finally, replace paragraph 1's text with the new content about removing
# the INTRODUCTION section. It is written as three runs so that the
# English heading "INTRODUCTION" does not carry the Hebrew/RTL run
# formatting used for the surrounding text (mirroring the source runs).
# ---------------------------------------------------------------------------
$p1final = $d.Paragraphs(1)
$r1final = $d.Range($p1final.Range.Start, $p1final.Range.End - 1)
$r1final.Text = "להעיף את ה- INTRODUCTION (לבדוק שחלקים ממנה מתחילים את הפרקים הרלוונטיים)"

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
